# Implements increasing the req_energy to range: replace row 26 with new data
# and append new rows 27-37 with the same column layout (B..X).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newData = @(
    @(213.4049414013323, 11.144, 13.40483300259625, 231246, 293.6298755186722, 239880, 28747.96, 163, 52, 4, 40, 49, 3, 1485, 0.28, 2.325, 0.015, 28000, 315, 400, 240, 90000, 50000),
    @(484.323441769288, 84, 4.036028681410734, 120600, 306.5161307454832, 50320.08, 58103.28, 111, 47, 14, 165, 47, 32, 1171, 0.28, 2.33, 0.015, 27700, 315, 398, 240, 90000, 50000),
    @(242.2955009345596, 11.144, 15.21956664161806, 227682, 294.4796680497925, 236910, 28923.92, 163, 53, 4, 40, 48, 3, 1171, 0.28, 2.33, 0.015, 27700, 315, 398, 240, 90000, 50000),
    @(277.1870672686526, 22.4, 8.662095852145395, 103635, 306.9875518672199, 122565, 33260, 163, 56, 4, 53, 38, 4, 1171, 0.28, 2.33, 0.015, 27700, 315, 398, 240, 90000, 50000),
    @(484.323441769288, 84, 4.036028681410734, 120600, 306.5161307454832, 50320.08, 58103.28, 111, 47, 14, 165, 47, 32, 1171, 0.28, 2.33, 0.015, 27700, 315, 398, 240, 90000, 50000),
    @(500.8068765952128, 84, 4.173390638293441, 180900, 308.6595079373123, 50404.68000000001, 60133.68, 111, 47, 20, 165, 47, 22, 1171, 0.28, 2.33, 0.015, 40000, 315, 398, 240, 90000, 50000),
    @(422.0531882892127, 42, 7.034219804820212, 113022, 283.35136806307, 50150.88, 50150.88, 165, 94, 26, 0, 0, 0, 1171, 0.28, 2.33, 0.015, 40000, 315, 398, 240, 90000, 50000),
    @(504.6793708753279, 84, 4.205661423961065, 228528, 309.0990589655484, 54648, 60609.60000000001, 111, 48, 21, 165, 46, 20, 1171, 0.28, 2.33, 0.015, 50000, 315, 398, 240, 90000, 50000),
    @(446.2492591759789, 42, 7.437487652932981, 140616, 301.9013757922399, 53434.08, 53434.07999999999, 165, 93, 28, 0, 0, 0, 1171, 0.28, 2.33, 0.015, 50000, 315, 398, 240, 90000, 50000),
    @(558.8017094888542, 126, 3.104453941604746, 217721.52, 309.2414937759336, 53777.52, 67113.36000000002, 111, 45, 22, 320, 48, 42, 1171, 0.28, 2.33, 0.015, 60000, 315, 398, 240, 90000, 50000),
    @(420.0390040416348, 50.40000000000001, 5.833875056133818, 171736, 308.2107883817427, 66866.5, 50426.3, 111, 47, 21, 96, 47, 10, 1171, 0.28, 2.33, 0.015, 27700, 315, 398, 240, 90000, 50000),
    @(549.0243392415929, 126, 3.050135218008849, 241544.16, 309.0406639004149, 54747.36000000001, 65933.64000000001, 111, 47, 25, 320, 46, 29, 1171, 0.28, 2.33, 0.015, 40000, 315, 398, 240, 90000, 50000),
)

$startRow = 26
for ($i = 0; $i -lt $newData.Count; $i++) {
    $rowVals = $newData[$i]
    $r = $startRow + $i
    for ($j = 0; $j -lt $rowVals.Count; $j++) {
        $ws.Cells.Item($r, 2 + $j).Value = $rowVals[$j]
    }
}

# Match the recorded selection state after the edit: rows 26:27 selected,
# active cell A26, sqref spans the full rows (A26:XFD27).
$ws.Range("A26:A27").EntireRow.Select()
